$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update feature names (column A) and importance values (column B)
# to reflect the re-ranked feature importances.

$ws.Range("A2").Value = "MACD"
$ws.Range("B2").Value = 0.4229088338509427

$ws.Range("A3").Value = "RSI"
$ws.Range("B3").Value = 0.2848388270457197

$ws.Range("A4").Value = "Signal_line"
$ws.Range("B4").Value = 0.1417336753192148

$ws.Range("A5").Value = "close_long"
$ws.Range("B5").Value = 0.03846227381309363

$ws.Range("A6").Value = "close_short"
$ws.Range("B6").Value = 0.02833904293109818

$ws.Range("A7").Value = "DJI"
$ws.Range("B7").Value = 0.02684635520520896

$ws.Range("A8").Value = "VIX_short"
$ws.Range("B8").Value = 0.02202210527213354

$ws.Range("A9").Value = "VIX_long"
$ws.Range("B9").Value = 0.01791150646154515

$ws.Range("A10").Value = "VIX"
$ws.Range("B10").Value = 0.01693738010104325
